# update master barang multi satuan
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# The "1234KDFH" placeholder barcode values in K3:K5 are no longer valid
# for these rows - clear them out (K2 keeps its value).
$ws.Range("K3").ClearContents()
$ws.Range("K4").ClearContents()
$ws.Range("K5").ClearContents()

# Scroll the view over a bit and move the selection to R8, matching the
# author's last on-screen position when the workbook was saved.
$excel.ActiveWindow.ScrollColumn = 2
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("R8").Select()
